$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Y"
$ws.Range("C2").Value = "Y"
$ws.Range("B3").Value = "Y"
$ws.Range("C3").Value = "Y"
$ws.Range("B4").Value = "Y"
$ws.Range("C4").Value = "Y"

$ws.Range("C2").Select() | Out-Null
